# Use paragraph symbol (¶) to separate adoc table rows.
# - Existing table rule descriptions (table1, table2, table5, table6) have the
#   newline that separated adoc table data-rows replaced with a pilcrow (¶).
# - table5's rule description is extended with 5 more rows (Name10..Name14)
#   and its trailing "..." placeholder is replaced with a proper "===" fence.
# - A brand new rule "table7" is inserted as a table (after table6), which
#   shifts every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$PILCROW = [string][char]0x00B6

# Leading apostrophe forces Excel to treat the assigned text as a literal
# string instead of trying to parse a leading "=" as the start of a formula
# (several of these rule-description strings legitimately start with "===").

# --- 1. table1 rule (row 29): "===`nWITH anchor`nWITHOUT anchor`n===" ---
$ws.Range("C29").Value = "'===`nWITH anchor" + $PILCROW + "WITHOUT anchor`n==="

# --- 2. table2 rule (row 30) ---
$ws.Range("C30").Value = "'Header 1|Header 2`n===`nCell in column 1, row 1|Cell in column 2, row 1" + $PILCROW + "Cell in column 1, row 2|Cell in column 2, row 2`n==="

# --- 3. table5 rule (row 33): extend with Name10..Name14 and join rows with pilcrow ---
$table5Rows = @("Roses|Red","Violets|Blue")
for ($i = 1; $i -le 14; $i++) {
    $table5Rows += "Name$i|Color$i"
}
$table5Body = [string]::Join($PILCROW, $table5Rows)
$ws.Range("C33").Value = "'Name|Color`n===`n" + $table5Body + "`n==="

# --- 4. table6 rule (row 34) ---
$ws.Range("C34").Value = "'X1|X2`n===`nA|B" + $PILCROW + "C|D`n==="

# --- 5. Insert new row 35 for the new "table7" rule, pushing everything else down ---
$ws.Rows.Item(35).Insert()

$ws.Range("A35").Value = "'my-chapter_name"
$ws.Range("B35").Value = "'table7"
$ws.Range("C35").Value = "'ColA|ColB`n===`n0`n1|Off`nOn`n==="
$ws.Range("D35").Value = "'[""norm:table:column-first-order""]"

# --- 6. Resize the worksheet table to include the new row ---
$lo = $ws.ListObjects.Item(1)
$loRange = $lo.Range
$lastRow = $loRange.Rows.Count + $loRange.Row - 1
$lastCol = $loRange.Columns.Count
$newLastRow = $lastRow + 1
$newRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($newLastRow, $lastCol))
$lo.Resize($newRange)
